# Update "想去人数" (want-to-go count) figures in column F, refreshed data
# from the gh-pages generator run.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 13188
$ws1.Range("F5").Value = 7
$ws1.Range("F6").Value = 105
$ws1.Range("F11").Value = 13136
$ws1.Range("F13").Value = 566
$ws1.Range("F14").Value = 8812
$ws1.Range("F15").Value = 7885
$ws1.Range("F21").Value = 5
$ws1.Range("F25").Value = 388
$ws1.Range("F27").Value = 78
$ws1.Range("F28").Value = 347

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 13188
$ws4.Range("F6").Value = 7
$ws4.Range("F7").Value = 105
$ws4.Range("F12").Value = 13136
$ws4.Range("F14").Value = 566
$ws4.Range("F15").Value = 8812
$ws4.Range("F16").Value = 7885
$ws4.Range("F22").Value = 5
$ws4.Range("F28").Value = 388
$ws4.Range("F30").Value = 78
$ws4.Range("F31").Value = 347
